$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check on" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 10.02.2022 01:15"

# Convert D5 from text "-0.4" to a real number
$ws.Range("D5").Value = -0.4

# Convert E5 from text timestamp to a real Excel date serial, matching the
# date formatting used by the other rows in column E (style used by E2:E4/E6:E10)
$ws.Range("E5").Value = 44602.04287037037
$ws.Range("E5").NumberFormat = "YYYY-MM-DD HH:MM:SS"
